$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (single decimal point, pure numeric text).
$textCells = @("D5", "D10", "D11", "D16", "D19", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D33", "D35", "D37", "D40", "D41", "D43", "D44", "D46", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.722.27"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.601.06"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "211.82"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "19.73"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "1.826.59"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "1.590.43"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "64.93"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "26.689.60"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "0.0" + ([string][char]0x2083) + "0740"
$ws.Range("D19").Value = "210.46"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("D24").Value = "8.95"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "144.03"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "7.07"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("D29").Value = "15.36"
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").Value = "0.0510"
$ws.Range("E30").Value = "  -0.34%  "
$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").Value = "2.97"
$ws.Range("D34").Value = "1.291.31"
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").Value = "0.599"
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("E38").Value = "  +5.94%  "
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").Value = "0.829"
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("D41").Value = "5.39"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "0.782"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").Value = "63.03"
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("D45").Value = "1.739.10"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("D46").Value = "90.61"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").Value = "0.0515"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").Value = "7.44"
$ws.Range("E51").Value = "  +0.31%  "
